$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2222
    $ws.Range("F3").Value = 1688
    $ws.Range("F5").Value = 1079
    $ws.Range("F6").Value = 752
    $ws.Range("F8").Value = 5799
    $ws.Range("F9").Value = 86
}
